# Fruta / hortaliza, semanal
#
# Insert two new daily price records at rows 621-622 of the "Hortaliza,
# Vega Central Mapocho de Santiago - Zapallo italiano" sheet. All existing
# rows from 621 downward shift down by two (to 623-660), preserving their
# original values and the date-format style on column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 621:658 down to 623:660, carrying formatting along.
$ws.Rows("621:622").Insert()

# New row 621
$ws.Range("A621").Value = 9
$ws.Range("B621").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C621").Value = "Metropolitana"
$ws.Range("D621").Value = 45267
$ws.Range("E621").Value = 13
$ws.Range("F621").Value = 100112032
$ws.Range("G621").Value = "Zapallo italiano"
$ws.Range("H621").Value = "Sin especificar"
$ws.Range("I621").Value = "Primera"
$ws.Range("J621").Value = 160
$ws.Range("K621").Value = 9000
$ws.Range("L621").Value = 10000
$ws.Range("M621").Value = 9500
$ws.Range("N621").Value = "`$/caja 50 unidades"
$ws.Range("O621").Value = "Región Metropolitana"
$ws.Range("P621").Value = 190
$ws.Range("Q621").Value = 50
$ws.Range("R621").Value = "Hortaliza"

# New row 622
$ws.Range("A622").Value = 9
$ws.Range("B622").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C622").Value = "Metropolitana"
$ws.Range("D622").Value = 45267
$ws.Range("E622").Value = 13
$ws.Range("F622").Value = 100112032
$ws.Range("G622").Value = "Zapallo italiano"
$ws.Range("H622").Value = "Sin especificar"
$ws.Range("I622").Value = "Primera"
$ws.Range("J622").Value = 250
$ws.Range("K622").Value = 9000
$ws.Range("L622").Value = 10000
$ws.Range("M622").Value = 9500
$ws.Range("N622").Value = "`$/caja 50 unidades"
$ws.Range("O622").Value = "Región de O'Higgins"
$ws.Range("P622").Value = 190
$ws.Range("Q622").Value = 50
$ws.Range("R622").Value = "Hortaliza"
